# Update quickstart_annotated examples to use the new WorkbookLocation()
# xlSlim function instead of the old CELL("filename",...) array-formula
# trick, and tidy up the dependent formula + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 used to derive the workbook folder with a volatile array formula
# built on CELL("filename",...)/FIND(...). Replace it with the simpler
# xlSlim-provided _xll.WorkbookLocation() array formula.
$ws.Range("B1").FormulaArray = "_xll.WorkbookLocation()"

# B2 builds the full path to the .py module from B1. WorkbookLocation()
# no longer returns a trailing backslash, so the separator is now added
# explicitly here.
$ws.Range("B2").Formula = '=B1&"\quickstart_annotated.py"'

# Leave the selection on B3 (the "register module" cell), matching the
# state the workbook was left in.
$ws.Range("B3").Select() | Out-Null
